$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Bmp2"
$ws.Cells.Item(2, 3).Value = "Bmpr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 9.163165666666666
$ws.Cells.Item(2, 8).Value = 27.489497
$ws.Cells.Item(2, 9).Value = 0.2800251397703982
$ws.Cells.Item(2, 10).Value = 0.2800251397703982
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 28.86405866666667
$ws.Cells.Item(2, 14).Value = 86.59217600000001
$ws.Cells.Item(2, 15).Value = 0.3025739760541936
$ws.Cells.Item(2, 16).Value = 0.3025739760541936
$ws.Cells.Item(2, 17).Value = 264.4861513750524
$ws.Cells.Item(2, 18).Value = 2380.375362375472
$ws.Cells.Item(2, 19).Value = 0.08472831993546066
$ws.Cells.Item(2, 20).Value = 0.08472831993546068

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Bmp2"
$ws.Cells.Item(3, 3).Value = "Bmpr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 9.163165666666666
$ws.Cells.Item(3, 8).Value = 27.489497
$ws.Cells.Item(3, 9).Value = 0.2800251397703982
$ws.Cells.Item(3, 10).Value = 0.2800251397703982
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 42.68037399999999
$ws.Cells.Item(3, 14).Value = 128.041122
$ws.Cells.Item(3, 15).Value = 0.4474066038250392
$ws.Cells.Item(3, 16).Value = 0.4474066038250392
$ws.Cells.Item(3, 17).Value = 391.0873376772926
$ws.Cells.Item(3, 18).Value = 3519.786039095633
$ws.Cells.Item(3, 19).Value = 0.1252850967703057
$ws.Cells.Item(3, 20).Value = 0.1252850967703058

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Bmp2"
$ws.Cells.Item(4, 3).Value = "Bmpr2"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 9.163165666666666
$ws.Cells.Item(4, 8).Value = 27.489497
$ws.Cells.Item(4, 9).Value = 0.2800251397703982
$ws.Cells.Item(4, 10).Value = 0.2800251397703982
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 23.85061433333334
$ws.Cells.Item(4, 14).Value = 71.551843
$ws.Cells.Item(4, 15).Value = 0.2500194201207672
$ws.Cells.Item(4, 16).Value = 0.2500194201207672
$ws.Cells.Item(4, 17).Value = 218.5471303881079
$ws.Cells.Item(4, 18).Value = 1966.924173492971
$ws.Cells.Item(4, 19).Value = 0.07001172306463174
$ws.Cells.Item(4, 20).Value = 0.07001172306463176

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Bmp2"
$ws.Cells.Item(5, 3).Value = "Bmpr2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 12.06704733333333
$ws.Cells.Item(5, 8).Value = 36.201142
$ws.Cells.Item(5, 9).Value = 0.3687673822623249
$ws.Cells.Item(5, 10).Value = 0.3687673822623249
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 28.86405866666667
$ws.Cells.Item(5, 14).Value = 86.59217600000001
$ws.Cells.Item(5, 15).Value = 0.3025739760541936
$ws.Cells.Item(5, 16).Value = 0.3025739760541936
$ws.Cells.Item(5, 17).Value = 348.3039621627769
$ws.Cells.Item(5, 18).Value = 3134.735659464992
$ws.Cells.Item(5, 19).Value = 0.1115794130902083
$ws.Cells.Item(5, 20).Value = 0.1115794130902083

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Bmp2"
$ws.Cells.Item(6, 3).Value = "Bmpr2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 12.06704733333333
$ws.Cells.Item(6, 8).Value = 36.201142
$ws.Cells.Item(6, 9).Value = 0.3687673822623249
$ws.Cells.Item(6, 10).Value = 0.3687673822623249
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 42.68037399999999
$ws.Cells.Item(6, 14).Value = 128.041122
$ws.Cells.Item(6, 15).Value = 0.4474066038250392
$ws.Cells.Item(6, 16).Value = 0.4474066038250392
$ws.Cells.Item(6, 17).Value = 515.0260932623693
$ws.Cells.Item(6, 18).Value = 4635.234839361323
$ws.Cells.Item(6, 19).Value = 0.1649889620994368
$ws.Cells.Item(6, 20).Value = 0.1649889620994368

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Bmp2"
$ws.Cells.Item(7, 3).Value = "Bmpr2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 12.06704733333333
$ws.Cells.Item(7, 8).Value = 36.201142
$ws.Cells.Item(7, 9).Value = 0.3687673822623249
$ws.Cells.Item(7, 10).Value = 0.3687673822623249
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 23.85061433333334
$ws.Cells.Item(7, 14).Value = 71.551843
$ws.Cells.Item(7, 15).Value = 0.2500194201207672
$ws.Cells.Item(7, 16).Value = 0.2500194201207672
$ws.Cells.Item(7, 17).Value = 287.8064920894118
$ws.Cells.Item(7, 18).Value = 2590.258428804706
$ws.Cells.Item(7, 19).Value = 0.09219900707267976
$ws.Cells.Item(7, 20).Value = 0.09219900707267976

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Bmp2"
$ws.Cells.Item(8, 3).Value = "Bmpr2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 11.49244066666667
$ws.Cells.Item(8, 8).Value = 34.477322
$ws.Cells.Item(8, 9).Value = 0.3512074779672769
$ws.Cells.Item(8, 10).Value = 0.351207477967277
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 28.86405866666667
$ws.Cells.Item(8, 14).Value = 86.59217600000001
$ws.Cells.Item(8, 15).Value = 0.3025739760541936
$ws.Cells.Item(8, 16).Value = 0.3025739760541936
$ws.Cells.Item(8, 17).Value = 331.7184816258525
$ws.Cells.Item(8, 18).Value = 2985.466334632672
$ws.Cells.Item(8, 19).Value = 0.1062662430285246
$ws.Cells.Item(8, 20).Value = 0.1062662430285246

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Bmp2"
$ws.Cells.Item(9, 3).Value = "Bmpr2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 11.49244066666667
$ws.Cells.Item(9, 8).Value = 34.477322
$ws.Cells.Item(9, 9).Value = 0.3512074779672769
$ws.Cells.Item(9, 10).Value = 0.351207477967277
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 42.68037399999999
$ws.Cells.Item(9, 14).Value = 128.041122
$ws.Cells.Item(9, 15).Value = 0.4474066038250392
$ws.Cells.Item(9, 16).Value = 0.4474066038250392
$ws.Cells.Item(9, 17).Value = 490.5016658261426
$ws.Cells.Item(9, 18).Value = 4414.514992435284
$ws.Cells.Item(9, 19).Value = 0.1571325449552966
$ws.Cells.Item(9, 20).Value = 0.1571325449552967

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Bmp2"
$ws.Cells.Item(10, 3).Value = "Bmpr2"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 11.49244066666667
$ws.Cells.Item(10, 8).Value = 34.477322
$ws.Cells.Item(10, 9).Value = 0.3512074779672769
$ws.Cells.Item(10, 10).Value = 0.351207477967277
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 23.85061433333334
$ws.Cells.Item(10, 14).Value = 71.551843
$ws.Cells.Item(10, 15).Value = 0.2500194201207672
$ws.Cells.Item(10, 16).Value = 0.2500194201207672
$ws.Cells.Item(10, 17).Value = 274.1017700893829
$ws.Cells.Item(10, 18).Value = 2466.915930804446
$ws.Cells.Item(10, 19).Value = 0.0878086899834557
$ws.Cells.Item(10, 20).Value = 0.08780868998345571
